$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")

# hunk 0 - sheet ALC
$ws_ALC.Range("H15").Value = 45.53
$ws_ALC.Range("I15").Value = 45.53
$ws_ALC.Range("K15").Value = 136.59
$ws_ALC.Range("M15").Value = 32.41

# hunk 1 - sheet ALC
$ws_ALC.Range("H17").Value = 4879.9116
$ws_ALC.Range("J17").Value = 4970.242
$ws_ALC.Range("L17").Value = 14910.726
$ws_ALC.Range("N17").Value = -15246.726

# hunk 2 - sheet ALC
$ws_ALC.Range("H19").Value = 288.41666
$ws_ALC.Range("J19").Value = 270
$ws_ALC.Range("L19").Value = 270
$ws_ALC.Range("N19").Value = -620

# hunk 3 - sheet ALC
$ws_ALC.Range("H88").Value = 698.9167
$ws_ALC.Range("I88").Value = 447.16666
$ws_ALC.Range("J88").Value = 950.6667
$ws_ALC.Range("K88").Value = 447.16666
$ws_ALC.Range("L88").Value = 950.6667
$ws_ALC.Range("M88").Value = -41.16665999999998
$ws_ALC.Range("N88").Value = -1762.6667

# hunk 4 - sheet ALC
$ws_ALC.Range("H91").Value = 698.9167
$ws_ALC.Range("I91").Value = 447.16666
$ws_ALC.Range("J91").Value = 950.6667
$ws_ALC.Range("K91").Value = 447.16666
$ws_ALC.Range("L91").Value = 950.6667
$ws_ALC.Range("M91").Value = 956.83334
$ws_ALC.Range("N91").Value = -3758.6667

# hunk 5 - sheet ALC
$ws_ALC.Range("H112").Value = 3704796
$ws_ALC.Range("J112").Value = 3969367.2
$ws_ALC.Range("L112").Value = 11908101.6
$ws_ALC.Range("N112").Value = -11910317.6

# hunk 6 - sheet ALC
$ws_ALC.Range("H132").Value = 2564.4866
$ws_ALC.Range("I132").Value = 2594.0557
$ws_ALC.Range("J132").Value = 1500
$ws_ALC.Range("K132").Value = 7782.1671
$ws_ALC.Range("L132").Value = 4500
$ws_ALC.Range("M132").Value = -5252.1671
$ws_ALC.Range("N132").Value = -9560

# hunk 7 - sheet ALC
$ws_ALC.Range("H137").Value = 66071.71000000001
$ws_ALC.Range("I137").Value = 79555.69
$ws_ALC.Range("J137").Value = 3555.0908
$ws_ALC.Range("K137").Value = 238667.07
$ws_ALC.Range("L137").Value = 10665.2724
$ws_ALC.Range("M137").Value = -236117.07
$ws_ALC.Range("N137").Value = -15765.2724

# hunk 8 - sheet ALC
$ws_ALC.Range("H141").Value = 1828.8518
$ws_ALC.Range("I141").Value = 1530.1818
$ws_ALC.Range("K141").Value = 4590.5454
$ws_ALC.Range("M141").Value = 589.4546

# hunk 9 - sheet ARM
$ws_ARM.Range("H45").Value = 2595.0833
$ws_ARM.Range("I45").Value = 2841.6667
$ws_ARM.Range("J45").Value = 2348.5
$ws_ARM.Range("K45").Value = 2841.6667
$ws_ARM.Range("L45").Value = 2348.5
$ws_ARM.Range("M45").Value = -2464.6667
$ws_ARM.Range("N45").Value = -3102.5

# hunk 10 - sheet ARM
$ws_ARM.Range("H61").Value = 2018.4103
$ws_ARM.Range("I61").Value = 1679.6177
$ws_ARM.Range("J61").Value = 4322.2
$ws_ARM.Range("K61").Value = 1679.6177
$ws_ARM.Range("L61").Value = 4322.2
$ws_ARM.Range("M61").Value = -1467.6177
$ws_ARM.Range("N61").Value = -4746.2

# hunk 11 - sheet ARM
$ws_ARM.Range("H97").Value = 554.6818
$ws_ARM.Range("I97").Value = 573.8946999999999
$ws_ARM.Range("J97").Value = 433
$ws_ARM.Range("K97").Value = 573.8946999999999
$ws_ARM.Range("L97").Value = 433
$ws_ARM.Range("M97").Value = -77.89469999999994
$ws_ARM.Range("N97").Value = -1425

# hunk 12 - sheet ARM
$ws_ARM.Range("H110").Value = 1115.6428
$ws_ARM.Range("I110").Value = 1073.2727
$ws_ARM.Range("K110").Value = 1073.2727
$ws_ARM.Range("M110").Value = 971.7273

# hunk 13 - sheet ARM
$ws_ARM.Range("H132").Value = 13338.934
$ws_ARM.Range("I132").Value = 1997.0322
$ws_ARM.Range("J132").Value = 38453.145
$ws_ARM.Range("K132").Value = 5991.096600000001
$ws_ARM.Range("L132").Value = 115359.435
$ws_ARM.Range("M132").Value = -3461.096600000001
$ws_ARM.Range("N132").Value = -120419.435

# hunk 14 - sheet ARM
$ws_ARM.Range("H136").Value = 2018.4103
$ws_ARM.Range("I136").Value = 1679.6177
$ws_ARM.Range("J136").Value = 4322.2
$ws_ARM.Range("K136").Value = 5038.8531
$ws_ARM.Range("L136").Value = 12966.6
$ws_ARM.Range("M136").Value = -2488.8531
$ws_ARM.Range("N136").Value = -18066.6

# hunk 15 - sheet BSM
$ws_BSM.Range("H20").Value = 2589.2856
$ws_BSM.Range("I20").Value = 2205.8823
$ws_BSM.Range("K20").Value = 2205.8823
$ws_BSM.Range("M20").Value = -1958.8823

# hunk 16 - sheet BSM
$ws_BSM.Range("H94").Value = 928.52
$ws_BSM.Range("I94").Value = 681.2
$ws_BSM.Range("K94").Value = 681.2
$ws_BSM.Range("M94").Value = -230.2

# hunk 17 - sheet CRP
$ws_CRP.Range("H7").Value = 43.714287
$ws_CRP.Range("I7").Value = 10
$ws_CRP.Range("J7").Value = 49.333332
$ws_CRP.Range("K7").Value = 10
$ws_CRP.Range("L7").Value = 49.333332
$ws_CRP.Range("M7").Value = 103
$ws_CRP.Range("N7").Value = -275.333332

# hunk 18 - sheet CRP
$ws_CRP.Range("H22").Value = 171
$ws_CRP.Range("I22").Value = 158.92308
$ws_CRP.Range("J22").Value = 186.7
$ws_CRP.Range("K22").Value = 158.92308
$ws_CRP.Range("L22").Value = 186.7
$ws_CRP.Range("M22").Value = 191.07692
$ws_CRP.Range("N22").Value = -886.7

# hunk 19 - sheet CRP
$ws_CRP.Range("H58").Value = 16537.182
$ws_CRP.Range("I58").Value = 1388.1111
$ws_CRP.Range("J58").Value = 34716.066
$ws_CRP.Range("K58").Value = 1388.1111
$ws_CRP.Range("L58").Value = 34716.066
$ws_CRP.Range("M58").Value = -1185.1111
$ws_CRP.Range("N58").Value = -35122.066

# hunk 20 - sheet CRP
$ws_CRP.Range("H94").Value = 4530.3335
$ws_CRP.Range("J94").Value = 6780
$ws_CRP.Range("L94").Value = 6780
$ws_CRP.Range("N94").Value = -7682

# hunk 21 - sheet CRP
$ws_CRP.Range("H105").Value = 856.125
$ws_CRP.Range("I105").Value = 833.3333
$ws_CRP.Range("J105").Value = 924.5
$ws_CRP.Range("K105").Value = 833.3333
$ws_CRP.Range("L105").Value = 924.5
$ws_CRP.Range("M105").Value = 913.6667
$ws_CRP.Range("N105").Value = -4418.5

# hunk 22 - sheet CRP
$ws_CRP.Range("H132").Value = 2793.037
$ws_CRP.Range("I132").Value = 2019.579
$ws_CRP.Range("K132").Value = 6058.737
$ws_CRP.Range("M132").Value = -3528.737

# hunk 23 - sheet CRP
$ws_CRP.Range("H136").Value = 16537.182
$ws_CRP.Range("I136").Value = 1388.1111
$ws_CRP.Range("J136").Value = 34716.066
$ws_CRP.Range("K136").Value = 4164.3333
$ws_CRP.Range("L136").Value = 104148.198
$ws_CRP.Range("M136").Value = -1614.3333
$ws_CRP.Range("N136").Value = -109248.198

# hunk 24 - sheet CUL
$ws_CUL.Range("H5").Value = 1150.9736
$ws_CUL.Range("J5").Value = 2094.8
$ws_CUL.Range("L5").Value = 6284.400000000001
$ws_CUL.Range("N5").Value = -6508.400000000001

# hunk 25 - sheet CUL
$ws_CUL.Range("H100").Value = 3102.3333
$ws_CUL.Range("J100").Value = 3102.3333
$ws_CUL.Range("L100").Value = 9306.999899999999
$ws_CUL.Range("N100").Value = -10928.9999

# hunk 26 - sheet CUL
$ws_CUL.Range("H131").Value = 651.9899
$ws_CUL.Range("J131").Value = 771.0273999999999
$ws_CUL.Range("L131").Value = 2313.0822
$ws_CUL.Range("N131").Value = -12393.0822

# hunk 27 - sheet CUL
$ws_CUL.Range("H135").Value = 1150.9736
$ws_CUL.Range("J135").Value = 2094.8
$ws_CUL.Range("L135").Value = 18853.2
$ws_CUL.Range("N135").Value = -23923.2

# hunk 28 - sheet GSM
$ws_GSM.Range("H102").Value = 3090.8235
$ws_GSM.Range("I102").Value = 2702
$ws_GSM.Range("J102").Value = 6007
$ws_GSM.Range("K102").Value = 2702
$ws_GSM.Range("L102").Value = 6007
$ws_GSM.Range("M102").Value = -1080
$ws_GSM.Range("N102").Value = -9251

# hunk 29 - sheet GSM
$ws_GSM.Range("H122").Value = 5241
$ws_GSM.Range("I122").Value = 4929.231
$ws_GSM.Range("J122").Value = 5916.5
$ws_GSM.Range("K122").Value = 14787.693
$ws_GSM.Range("L122").Value = 17749.5
$ws_GSM.Range("M122").Value = -12337.693
$ws_GSM.Range("N122").Value = -22649.5

# hunk 30 - sheet LTW
$ws_LTW.Range("H132").Value = 187304.03
$ws_LTW.Range("I132").Value = 242411.16
$ws_LTW.Range("J132").Value = 3613.6
$ws_LTW.Range("K132").Value = 727233.48
$ws_LTW.Range("L132").Value = 10840.8
$ws_LTW.Range("M132").Value = -724703.48
$ws_LTW.Range("N132").Value = -15900.8
